# Rerun analyses for misinfotext: add a "2015" sheet (between 2014 and 2016)
# and a corresponding summary row in the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2015" worksheet positioned right before "2016".
# ---------------------------------------------------------------------------
$sheet2016 = $wb.Worksheets.Item("2016")
$newSheet = $wb.Worksheets.Add($sheet2016)
$newSheet.Name = "2015"

# Header row, styled like the header rows on the other year sheets
# (bold font, thin border all around, centered / top aligned).
$newSheet.Range("A1").Value = "factcheckURL"
$newSheet.Range("B1").Value = "lexical_diversity"
$newSheet.Range("A1:B1").Font.Bold = $true
$newSheet.Range("A1:B1").Borders.LineStyle = 1
$newSheet.Range("A1:B1").HorizontalAlignment = -4108
$newSheet.Range("A1:B1").VerticalAlignment = -4160

$urls = @(
    "https://www.politifact.com/factchecks/2015/nov/06/greg-abbott/greg-abbott-embarrassed-says-californians-buying-m/",
    "https://www.politifact.com/factchecks/2015/oct/16/scott-walker/gov-scott-walker-says-board-wanted-accept-mickey-m/",
    "https://www.politifact.com/factchecks/2015/oct/05/dana-loesch/Planned-parenthood-86-percent-abortion-revenue/",
    "https://www.politifact.com/factchecks/2015/sep/23/donald-trump/hillary-clinton-obama-birther-fact-check/",
    "https://www.politifact.com/factchecks/2015/sep/10/ted-cruz/ted-cruz-says-deal-will-facilitate-and-accelerate-/",
    "https://www.politifact.com/factchecks/2015/aug/26/hillary-clinton/hillary-clinton-says-no-gop-candidate-has-talked-a/",
    "https://www.politifact.com/factchecks/2015/jul/21/wisconsin-state-afl-cio/wisconsin-afl-cio-says-scott-walker-budget-means-n/",
    "https://www.politifact.com/factchecks/2015/mar/20/glenn-beck/glenn-beck-says-barack-obama-took-iran-hamas-us-te/",
    "https://www.politifact.com/factchecks/2015/nov/08/ben-carson/ben-carson-said-no-one-who-signed-declaration-inde/",
    "https://www.politifact.com/factchecks/2015/jul/07/dinesh-dsouza/hillary-clinton-confederate-battle-flag-nope-old-i/",
    "https://www.politifact.com/factchecks/2015/jun/25/gavin-mcinnes/tweet-civil-war-was-about-secession-not-slavery/",
    "https://www.politifact.com/factchecks/2015/feb/06/scott-walker/despite-deliberate-actions-scott-walker-calls-chan/"
)

$values = @(
    1,
    1,
    0.9375,
    0.875,
    0.9333333333333333,
    0.8636363636363636,
    0.92,
    0.421602787456446,
    0.4005305039787798,
    1,
    0.875,
    0.8181818181818182
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $urlCell = $newSheet.Range("A$row")
    $urlCell.Value = $urls[$i]
    $newSheet.Hyperlinks.Add($urlCell, $urls[$i])
    $urlCell.Style = "Hyperlink"
    $newSheet.Range("B$row").Value = $values[$i]
}

# ---------------------------------------------------------------------------
# 2. Insert a "2015" row into the "Summary" sheet, right before "2016".
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Rows.Item(6).Insert()

$summary.Range("A6").Value = 2015
$summary.Range("A6").Font.Bold = $true

$summary.Range("B6").Value = 12
$summary.Range("C6").Value = 0.837065400548895
$summary.Range("D6").Value = 0.2073618696352768
$summary.Range("E6").Value = 0.4005305039787798
$summary.Range("F6").Value = 0.8522727272727273
$summary.Range("G6").Value = 0.8975
$summary.Range("H6").Value = 0.953125
$summary.Range("I6").Value = 1

Write-Output "Inserted 2015 sheet and Summary row."
